# Fix a typo in column K ("TERÇA FEIRA - 19H -  IECG CENTRO" with a double
# space) so it matches the already-existing single-space variant of the
# same text used elsewhere in the sheet. This removes the now-duplicate
# shared string from the workbook on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$correctedValue = "TERÇA FEIRA - 19H - IECG CENTRO"
$rowsToFix = @(2, 4, 8, 10, 15, 28, 36, 43, 55, 57, 70)
foreach ($r in $rowsToFix) {
    $ws.Cells.Item($r, 11).Value = $correctedValue
}

# Give K73/K74 the same (white fill) cell style already used by K70, so the
# last two rows visually match the rest of the filtered-in rows.
$ws.Range("K73").Style = $ws.Range("K70").Style
$ws.Range("K74").Style = $ws.Range("K70").Style

# Re-apply the AutoFilter over the full data range (it now reaches row 74)
# and filter column K ("DIA DE AULA") down to the corrected value only.
$fullRange = $ws.Range("A1:K74")
$fullRange.AutoFilter()
$fullRange.AutoFilter()
$criteria = @($correctedValue)
$fullRange.AutoFilter(11, $criteria, 7)
